$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, shifting the existing rows 133-222 down to 134-223
$ws.Rows(133).Insert()

# Populate the newly inserted row 133 with the new weekly price record
$ws.Range("A133").Value = 11
$ws.Range("B133").Value = "Vega Monumental Concepción"
$ws.Range("C133").Value = "Bíobío"
$ws.Range("D133").Value = 44806
$ws.Range("E133").Value = 8
$ws.Range("F133").Value = 100112040
$ws.Range("G133").Value = "Cilantro"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 220
$ws.Range("K133").Value = 5000
$ws.Range("L133").Value = 5500
$ws.Range("M133").Value = 5273
$ws.Range("N133").Value = "`$/caja 36 atados"
$ws.Range("O133").Value = "Región Metropolitana"
$ws.Range("P133").Value = 146
$ws.Range("Q133").Value = 36
$ws.Range("R133").Value = "Hortaliza"
